$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "71.946.21"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +3.86%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.616.46"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +6.89%  "

$ws.Range("E4").Value = "  -0.02%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "599.68"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.13%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "182.57"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.73%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "3.610.28"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +6.87%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.607"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.84%  "

$ws.Range("E9").Value = "  +0.08%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.206"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +6.22%  "

$ws.Range("E11").Value = "  +3.10%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "50.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.89%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000290"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.15%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "710.50"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.46%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.192.89"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +6.83%  "

$ws.Range("E16").Value = "  +3.88%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "72.091.57"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.584.46"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.69%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.123"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.79%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "18.57"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +5.12%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.56%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.936"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.65%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.79"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.79%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "17.72"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.45%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "105.58"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "

$ws.Range("E26").Value = "  +3.05%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.61%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.03"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.38%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "35.56"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.44%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.11"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +4.54%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.48"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +7.67%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.10"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +15.46%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "592.49"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.13%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "11.36"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.10%  "

$ws.Range("E35").Value = "  +1.53%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "60.12"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.55%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.657.05"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "

$ws.Range("E39").Value = "  +3.92%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0₃0785"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +12.25%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.09"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.81%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.45"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +5.08%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.80"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.08%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.0453"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +7.53%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.347"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +2.34%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.39"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.35%  "

$ws.Range("E47").Value = "  +4.12%  "

$ws.Range("E48").Value = "  +5.52%  "

$ws.Range("E49").Value = "  +2.01%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.14%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "133.12"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.37%  "
